$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New university rows appended to the data table (rows 10, 12-17).
# Columns: A=University B=Location C=Scholarship D=Fee Structure E=Sector
#          F=MBBS G=BDS H=D-Pharm I=Food Sciences J=Computer Science
#          K=Software Engineering L=Data Science/AI M=BS English N=BS Chemistry
#          O=BS Physics P=Accounting Q=BBA R=LLB S=Sport Activities

$rows = @(
    @{ Row = 10; Values = @("University of Health Sciences (UHS)", "Lahore", 0, 50000, "Public", 1, 1, 1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0) },
    @{ Row = 12; Values = @("Bahauddin Zakariya Univeristy", "Multan", 1, 40000, "Public", 0, 0, 1, 1, 1, 0, 0, 1, 1, 1, 1, 1, 1, 1) },
    @{ Row = 13; Values = @("UET", "Lahore", 0, 54000, "Public", 0, 0, 0, 0, 1, 1, 1, 0, 1, 1, 0, 1, 0, 0) },
    @{ Row = 14; Values = @("UET", "Faislabad", 0, 54000, "Public", 0, 0, 0, 0, 0, 0, 1, 0, 1, 1, 0, 1, 0, 0) },
    @{ Row = 15; Values = @("Comsats University", "Sahiwal", 1, 120000, "Public", 0, 0, 0, 1, 1, 1, 0, 0, 0, 0, 0, 1, 0, 1) },
    @{ Row = 16; Values = @("University of Sahiwal", "Sahiwal", 0, 30000, "Public", 0, 0, 0, 0, 1, 1, 0, 1, 1, 1, 0, 1, 1, 0) },
    @{ Row = 17; Values = @("NUST", "Islamabad", 1, 200000, "Public", 1, 0, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1) }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Widen column A to fit the newly added, longer university names.
$ws.Columns.Item(1).ColumnWidth = 33 - 5/6

# Leave the selection on the last populated cell, matching the final edit state.
$ws.Range("S17").Select() | Out-Null
